# Generate Report for Handoff
# Mark "b.md" as ready for handoff with a fresh handoff file / timestamp,
# on the Overview sheet as well as the per-locale (zh-cn / de-de) sheets.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: row 3 is b.md
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value2 = $statusReady
$wsOverview.Range("C3").Value2 = $statusReady
$wsOverview.Range("D3").Value2 = "2016-26-18 00:26:05"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is b.md
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhCnFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

$wsZhCn.Range("C3").Value2 = $statusReady
$wsZhCn.Range("D3").Value2 = $zhCnFile
$wsZhCn.Range("E3").Value2 = "2016-03-18 00:26:01"

foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = $zhCnFile
    }
}

# ---------------------------------------------------------------------
# de-de sheet: row 3 is b.md
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$deDeFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"

$wsDeDe.Range("C3").Value2 = $statusReady
$wsDeDe.Range("D3").Value2 = $deDeFile
$wsDeDe.Range("E3").Value2 = "2016-03-18 00:26:05"

foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = $deDeFile
    }
}
